$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing totals in row 9 (anhlavodich)
$ws.Range("B9").Value = 193100
$ws.Range("C9").Value = 6
$ws.Range("E9").Value = 6

# Update the accumulated history strings (shared strings content change)
$ws.Range("F9").Value = ";0;0;0;0;0;0"
$ws.Range("G9").Value = ";12;41;24;24;13;44"
$ws.Range("H9").Value = ";-500;-500;-500;-500;-1000;-1000"

# Add the new "vatpham" (item) history column
$ws.Range("I9").Value = ";1;2;3;4;2"
